# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 151
    3  = 132
    7  = 37
    8  = 502
    9  = 41
    10 = 1929
    12 = 75
    13 = 4031
    15 = 266
    17 = 72
    18 = 12
    20 = 2750
    21 = 38
    22 = 377
    25 = 57
    27 = 46
    30 = 36
    31 = 178
    32 = 156
    33 = 1583
    34 = 212
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
